# Fruta / hortaliza, semanal
# Insert two new weekly rows of data (Femacal de La Calera - Frutilla) above the
# existing row 138, shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 138; existing rows 138.. shift down to 140..
$ws.Rows("138:139").Insert()

# ---- Row 138 ----
$ws.Range("A138").Value = 3
$ws.Range("B138").Value = "Femacal de La Calera"
$ws.Range("C138").Value = "Coquimbo"
$ws.Range("D138").Value = 44582
$ws.Range("E138").Value = 5
$ws.Range("F138").Value = "Fruta"
$ws.Range("G138").Value = 100101
$ws.Range("H138").Value = "Berries"
$ws.Range("I138").Value = 100112025
$ws.Range("J138").Value = "Frutilla"
$ws.Range("K138").Value = "Sin especificar"
$ws.Range("L138").Value = "Primera"
$ws.Range("M138").Value = 230
$ws.Range("N138").Value = 5000
$ws.Range("O138").Value = 5500
$ws.Range("P138").Value = 5239
$ws.Range("Q138").Value = "$/bandeja 7 kilos"
$ws.Range("R138").Value = "Provincia de Melipilla"
$ws.Range("S138").Value = 748
$ws.Range("T138").Value = 7

# ---- Row 139 ----
$ws.Range("A139").Value = 3
$ws.Range("B139").Value = "Femacal de La Calera"
$ws.Range("C139").Value = "Coquimbo"
$ws.Range("D139").Value = 44582
$ws.Range("E139").Value = 5
$ws.Range("F139").Value = "Fruta"
$ws.Range("G139").Value = 100101
$ws.Range("H139").Value = "Berries"
$ws.Range("I139").Value = 100112025
$ws.Range("J139").Value = "Frutilla"
$ws.Range("K139").Value = "Sin especificar"
$ws.Range("L139").Value = "Segunda"
$ws.Range("M139").Value = 90
$ws.Range("N139").Value = 4000
$ws.Range("O139").Value = 4000
$ws.Range("P139").Value = 4000
$ws.Range("Q139").Value = "$/bandeja 7 kilos"
$ws.Range("R139").Value = "Provincia de Melipilla"
$ws.Range("S139").Value = 571
$ws.Range("T139").Value = 7
